$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.020572
$ws.Range("H2").Value = 0.061716
$ws.Range("I2").Value = 0.005721785283044617
$ws.Range("J2").Value = 0.005721785283044618
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.005969774965333334
$ws.Range("R2").Value = 0.053727974688
$ws.Range("S2").Value = 0.0001962223081794258
$ws.Range("T2").Value = 0.0001962223081794258

$ws.Range("G3").Value = 0.020572
$ws.Range("H3").Value = 0.061716
$ws.Range("I3").Value = 0.005721785283044617
$ws.Range("J3").Value = 0.005721785283044618
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 0.1462741339146666
$ws.Range("R3").Value = 1.316467205232
$ws.Range("S3").Value = 0.004807927995671053
$ws.Range("T3").Value = 0.004807927995671053

$ws.Range("G4").Value = 0.020572
$ws.Range("H4").Value = 0.061716
$ws.Range("I4").Value = 0.005721785283044617
$ws.Range("J4").Value = 0.005721785283044618
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 0.02183298816933333
$ws.Range("R4").Value = 0.196496893524
$ws.Range("S4").Value = 0.0007176349791941398
$ws.Range("T4").Value = 0.0007176349791941397

$ws.Range("H5").Value = 7.905868999999999
$ws.Range("I5").Value = 0.7329652747079958
$ws.Range("J5").Value = 0.7329652747079958
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 0.7647329515102222
$ws.Range("R5").Value = 6.882596563591999
$ws.Range("S5").Value = 0.02513623474211175
$ws.Range("T5").Value = 0.02513623474211175

$ws.Range("H6").Value = 7.905868999999999
$ws.Range("I6").Value = 0.7329652747079958
$ws.Range("J6").Value = 0.7329652747079958
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.6158994247068492
$ws.Range("T6").Value = 0.6158994247068492

$ws.Range("H7").Value = 7.905868999999999
$ws.Range("I7").Value = 0.7329652747079958
$ws.Range("J7").Value = 0.7329652747079958
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 2.796823260504556
$ws.Range("S7").Value = 0.09192961525903486
$ws.Range("T7").Value = 0.09192961525903483

$ws.Range("I8").Value = 0.2613129400089597
$ws.Range("J8").Value = 0.2613129400089597
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 0.2726385857235556
$ws.Range("R8").Value = 2.453747271512
$ws.Range("S8").Value = 0.008961438730959463
$ws.Range("T8").Value = 0.008961438730959463

$ws.Range("I9").Value = 0.2613129400089597
$ws.Range("J9").Value = 0.2613129400089597
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.2195772364305951
$ws.Range("T9").Value = 0.2195772364305951

$ws.Range("I10").Value = 0.2613129400089597
$ws.Range("J10").Value = 0.2613129400089597
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 0.9971087773278889
$ws.Range("R10").Value = 8.973978995951
$ws.Range("S10").Value = 0.03277426484740514
$ws.Range("T10").Value = 0.03277426484740513
